$d = $word.ActiveDocument

$replacements = @(
    @("90÷7=", "15÷4="),
    @("25÷8=", "49÷6="),
    @("97÷9=", "17÷5="),
    @("46÷4=", "85÷7="),
    @("28÷4=", "92÷4="),
    @("15÷3=", "32÷9="),
    @("33÷9=", "98÷2="),
    @("24÷9=", "94÷2="),
    @("40÷6=", "87÷7="),
    @("45÷3=", "25÷7="),
    @("15÷6=", "92÷2="),
    @("35÷5=", "83÷8="),
    @("92÷9=", "63÷5="),
    @("33÷8=", "47÷2="),
    @("11÷6=", "98÷8="),
    @("64÷2=", "84÷2="),
    @("65÷4=", "77÷3="),
    @("73÷6=", "34÷8="),
    @("39÷3=", "91÷7="),
    @("90÷4=", "44÷9="),
    @("22÷4=", "53÷9="),
    @("77÷7=", "14÷4="),
    @("15÷7=", "22÷8="),
    @("20÷9=", "21÷2="),
    @("44÷2=", "55÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
